# Applies updated Price (D) and Volume(1h) (E) values for the crypto symbol list.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "'325.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-2.17%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'44.46"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'0.84%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.491"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-6.18%"
$ws.Range("E4").Style = "Normal"
$ws.Range("E5").Value = "'-3.30%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'8.643"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-1.93%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.913"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-3.05%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'4.274"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-5.01%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.704"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-6.70%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.9402"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'0.37%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1172"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-6.80%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.1868"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-4.30%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.1009"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'4.49%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.04254"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'3.53%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.1065"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.11%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.001286"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-1.17%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.005902"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.99%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.586"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'2.33%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'-0.32%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'8.424"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-5.77%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1377"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'0.41%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.2526"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-1.80%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04239"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-3.71%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001234"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-1.89%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004562"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'3.35%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'-0.96%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0003988"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'-0.11%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.02637"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-6.54%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.05491"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-3.52%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007676"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-3.26%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1395"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-2.30%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.007166"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-20.74%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002030"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-3.40%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.009184"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-12.37%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00007102"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-2.22%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'-0.12%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.003526"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'8.96%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.002269"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-0.48%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002099"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.12%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0001999"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.12%"
$ws.Range("E51").Style = "Normal"
